$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 updates ---
$ws.Range("A2").Value = "098/TTTTT"
$ws.Range("H2").Value = 10000
$ws.Range("J2").Value = 1500
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 8500

# --- Row 3 is replaced with blank placeholder data (matches the "totals" row style) ---
$ws.Range("A3:F3").Value = " "
$ws.Range("G3").Value = " "
$ws.Range("H3").Value = 10000
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 1500
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 8500

# --- Old totals row 4 is removed entirely, shrinking the sheet to A1:M3 ---
$ws.Rows.Item(4).Delete()
